$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("H32").Value = 4
